$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.818.69'
$ws.Range('E2').Value = '  -1.25%  '

$ws.Range('D3').Value = '2.595.42'
$ws.Range('E3').Value = '  -1.96%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '551.09'
$ws.Range('E5').Value = '  +2.58%  '

$ws.Range('D6').Value = '143.13'
$ws.Range('E6').Value = '  -1.98%  '

$ws.Range('E7').Value = '  +0.00%  '

$ws.Range('E8').Value = '  +5.75%  '

$ws.Range('D9').Value = '6.78'
$ws.Range('E9').Value = '  +1.18%  '

$ws.Range('E10').Value = '  -1.75%  '

$ws.Range('E11').Value = '  +5.00%  '

$ws.Range('E12').Value = '  -0.74%  '

$ws.Range('D13').Value = '3.055.05'
$ws.Range('E13').Value = '  -1.76%  '

$ws.Range('D14').Value = '58.770.29'
$ws.Range('E14').Value = '  -1.14%  '

$ws.Range('D15').Value = '20.85'
$ws.Range('E15').Value = '  -2.14%  '

$ws.Range('D16').Value = '2.607.66'
$ws.Range('E16').Value = '  -0.35%  '

$ws.Range('E17').Value = '  -2.16%  '

$ws.Range('E18').Value = '  +1.26%  '

$ws.Range('D19').Value = '336.61'
$ws.Range('E19').Value = '  -1.05%  '

$ws.Range('E20').Value = '  -2.62%  '

$ws.Range('E21').Value = '  -1.61%  '

$ws.Range('E22').Value = '  +0.04%  '

$ws.Range('D23').Value = '66.77'
$ws.Range('E23').Value = '  +0.07%  '

$ws.Range('D24').Value = '0.428'
$ws.Range('E24').Value = '  +2.69%  '

$ws.Range('E25').Value = '  -0.03%  '

$ws.Range('E26').Value = '  -2.99%  '

$ws.Range('E27').Value = '  -2.01%  '

$ws.Range('E28').Value = '  +0.91%  '

$ws.Range('E29').Value = '  +0.00%  '

$ws.Range('E30').Value = '  +1.45%  '

$ws.Range('E31').Value = '  +2.25%  '

$ws.Range('D32').Value = '154.17'
$ws.Range('E32').Value = '  +1.84%  '

$ws.Range('D33').Value = '18.92'
$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('E34').Value = '  -2.05%  '

$ws.Range('D35').Value = '0.884'
$ws.Range('E35').Value = '  +4.25%  '

$ws.Range('E36').Value = '  -1.40%  '

$ws.Range('D37').Value = '37.01'
$ws.Range('E37').Value = '  -0.79%  '

$ws.Range('E38').Value = '  +0.89%  '

$ws.Range('D39').Value = '0.825'
$ws.Range('E39').Value = '  -1.45%  '

$ws.Range('E40').Value = '  +0.12%  '

$ws.Range('D41').Value = '283.77'
$ws.Range('E41').Value = '  -0.92%  '

$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  +0.08%  '

$ws.Range('E43').Value = '  -1.19%  '

$ws.Range('D44').Value = '0.0958'
$ws.Range('E44').Value = '  +1.34%  '

$ws.Range('D46').Value = '0.0533'
$ws.Range('E46').Value = '  -1.11%  '

$ws.Range('D47').Value = '0.0226'
$ws.Range('E47').Value = '  -0.12%  '

$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').Value = '1.941.64'
$ws.Range('E48').Value = '  -1.32%  '

$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '118.95'
$ws.Range('E49').Value = '  +7.06%  '

$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').Value = '17.83'
$ws.Range('E50').Value = '  -2.51%  '

$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '4.40'
$ws.Range('E51').Value = '  -3.74%  '
